$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Answer cycle repeats every 5 rows: B, B, D, C, A
$answers = @("B", "B", "D", "C", "A")

# Fill in question rows 27-35 (question numbers 26-34), mirroring the
# already-populated rows above (question #, answer letter, 2-point value).
for ($r = 27; $r -le 35; $r++) {
    $qNum = $r - 1
    $answer = $answers[($qNum - 1) % 5]
    $ws.Cells.Item($r, 1).Value = $qNum
    $ws.Cells.Item($r, 2).Value = $answer
    $ws.Cells.Item($r, 3).Value = 2
}

# Update the active selection to reflect the cell the author ended up on
# after entering the new answers (keeps the existing frozen header row).
$ws.Range("B31").Select()
